$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = 700
    4  = 800
    8  = 800
    10 = 800
    11 = 800
    12 = 800
    13 = 800
    14 = 800
    16 = 400
    17 = 500
    18 = 800
    19 = 800
    20 = 500
    22 = 366.67
    23 = 800
    25 = 800
    26 = 800
    29 = 800
    30 = 800
    31 = 800
    32 = 800
    34 = 800
    35 = 800
    37 = 800
    38 = 800
    39 = 800
    40 = 800
    41 = 800
    42 = 500
    43 = 800
    44 = 500
    45 = 800
    46 = 450
    47 = 650
    48 = 700
    49 = 750
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
